# The "Förändrad" (changed) date stamp in column C is bumped by one day
# (2023-09-19 -> 2023-09-20, i.e. Excel serial 45188 -> 45189) for every
# data row (rows 2-307) on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 307
$newValue = 45189

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newValue
}
